$d = $word.ActiveDocument

# 1. Merge "Number Plate " + "Detection" runs into a single run "Number Plate Detection"
$d.Content.Find.Execute("Number Plate Detection", $true, $false, $false, $false, $false, $true, 1, $false, "Number Plate Detection", 2) | Out-Null

# 2. Update the OS line from "Windows 10" to "macOS Monterey - Version 12.4" (en dash)
$enDash = [char]0x2013
$newOs = "macOS Monterey " + $enDash + " Version 12.4"
$d.Content.Find.Execute("Windows 10", $true, $false, $false, $false, $false, $true, 1, $false, $newOs, 2) | Out-Null

# 3. Remove the leftover "_GoBack" bookmark
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
